# Re-analyzed carbohydrate esters & amides:
# Ashish originally had 2 categories for carbohydrate esters and 2
# categories for amides (each occupying a specific FTIR spectral range).
# Those had been combined into a single "carboEster" / "amide" category;
# this change splits them back into their original categories
# ("carboEster1"/"carboEster2" and "amide1"/"amide2"), re-analyzed with
# factorial ANOVAs and Tukey's.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("litterChemistry")

# Current layout (before):
#  1 header
#  2 glycosidicBond
#  3 C_O_stretching
#  4 carboEster
#  5 lipid
#  6 alkane
#  7 amide
#
# Target layout (after):
#  1 header
#  2 glycosidicBond
#  3 C_O_stretching
#  4 carboEster
#  5 carboEster1
#  6 carboEster2
#  7 lipid
#  8 alkane
#  9 amide
# 10 amide1
# 11 amide2

# Insert the amide split first (bottom-up) so earlier row numbers don't
# shift while we still need them.
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).Insert()

# Insert two rows right after row 4 (carboEster), pushing
# lipid/alkane/amide.. down by two.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# Carry the formatting (font/border) of the row directly above into each
# freshly inserted row so the thin-border style is preserved exactly.
$ws.Range("A4:H4").Copy()
$ws.Range("A5:H5").PasteSpecial(-4122)
$ws.Range("A6:H6").PasteSpecial(-4122)

$ws.Range("A7:H7").Copy()
$ws.Range("A10:H10").PasteSpecial(-4122)
$ws.Range("A11:H11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 5 -> carboEster1
$ws.Range("A5").Value = "carboEster1"
$ws.Range("B5").Value = "o"
$ws.Range("C5").Value = "***"
$ws.Range("D5").Value = "o"
$ws.Range("E5").Value = "o"
$ws.Range("F5").Value = "o"
$ws.Range("G5").Value = "o"
$ws.Range("H5").Value = "o"

# Row 6 -> carboEster2
$ws.Range("A6").Value = "carboEster2"
$ws.Range("B6").Value = "o"
$ws.Range("C6").Value = "***"
$ws.Range("D6").Value = "***"
$ws.Range("E6").Value = "o"
$ws.Range("F6").Value = "o"
$ws.Range("G6").Value = "**"
$ws.Range("H6").Value = "o"

# Row 10 -> amide1
$ws.Range("A10").Value = "amide1"
$ws.Range("B10").Value = "*"
$ws.Range("C10").Value = "o"
$ws.Range("D10").Value = "***"
$ws.Range("E10").Value = "o"
$ws.Range("F10").Value = "o"
$ws.Range("G10").Value = "o"
$ws.Range("H10").Value = "o"

# Row 11 -> amide2
$ws.Range("A11").Value = "amide2"
$ws.Range("B11").Value = "o"
$ws.Range("C11").Value = "o"
$ws.Range("D11").Value = "***"
$ws.Range("E11").Value = "o"
$ws.Range("F11").Value = "o"
$ws.Range("G11").Value = "o"
$ws.Range("H11").Value = "o"

# The first cell of each new amide row had its fill explicitly toggled
# (cleared) while re-formatting, distinguishing it from the untouched
# rows above/below.
$ws.Range("A10").Interior.ColorIndex = -4142
$ws.Range("A11").Interior.ColorIndex = -4142

$ws.Range("H14").Select()
